$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, shifting rows 134:208 down to 135:209
$ws.Rows.Item(134).Insert()

# Populate the new row 134 with the inserted record's data.
# Static columns copied from neighboring rows.
$ws.Range("A134").Value = 9
$ws.Range("B134").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C134").Value = "Metropolitana"
$ws.Range("D134").Value = 44438
$ws.Range("D134").NumberFormat = $ws.Range("D135").NumberFormat
$ws.Range("E134").Value = 13
$ws.Range("F134").Value = 100112012
$ws.Range("G134").Value = "Espinaca"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 160
$ws.Range("K134").Value = 8000
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = 8500
$ws.Range("N134").Value = "`$/cuna 10 kilos"
$ws.Range("O134").Value = "Provincia de Chacabuco"
$ws.Range("P134").Value = 850
$ws.Range("Q134").Value = 10
$ws.Range("R134").Value = "Hortaliza"
